# Fruta / hortaliza, semanal
# Insert a new weekly record as row 105, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 105 (existing rows 105..181 shift to 106..182)
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new weekly observation
$ws.Range("A105").Value = 1
$ws.Range("B105").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C105").Value = "Arica y Parinacota"
$ws.Range("D105").Value = 44566
$ws.Range("E105").Value = 15
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100102
$ws.Range("H105").Value = "Cítricos"
$ws.Range("I105").Value = 100102003
$ws.Range("J105").Value = "Limón"
$ws.Range("K105").Value = "Sin especificar"
$ws.Range("L105").Value = "2a amarillo"
$ws.Range("M105").Value = 300
$ws.Range("N105").Value = 24000
$ws.Range("O105").Value = 25000
$ws.Range("P105").Value = 24500
$ws.Range("Q105").Value = "$/caja 20 kilos"
$ws.Range("R105").Value = "Región de Coquimbo"
$ws.Range("S105").Value = 1225
$ws.Range("T105").Value = 20

# Make sure the date column keeps the date number format used by the rest of column D
$ws.Range("D105").NumberFormat = $ws.Range("D106").NumberFormat
